$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-02-24 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-25 Sunday", 2) | Out-Null
$d.Content.Find.Execute("86×46=3956", $true, $false, $false, $false, $false, $true, 1, $false, "28×97=2716", 2) | Out-Null
$d.Content.Find.Execute("49×30=1470", $true, $false, $false, $false, $false, $true, 1, $false, "40×83=3320", 2) | Out-Null
$d.Content.Find.Execute("46×76=3496", $true, $false, $false, $false, $false, $true, 1, $false, "66×74=4884", 2) | Out-Null
$d.Content.Find.Execute("25×44=1100", $true, $false, $false, $false, $false, $true, 1, $false, "89×25=2225", 2) | Out-Null
$d.Content.Find.Execute("65×61=3965", $true, $false, $false, $false, $false, $true, 1, $false, "89×48=4272", 2) | Out-Null
$d.Content.Find.Execute("76×86=6536", $true, $false, $false, $false, $false, $true, 1, $false, "81×75=6075", 2) | Out-Null
$d.Content.Find.Execute("57×34=1938", $true, $false, $false, $false, $false, $true, 1, $false, "65×26=1690", 2) | Out-Null
$d.Content.Find.Execute("15×60=900", $true, $false, $false, $false, $false, $true, 1, $false, "40×62=2480", 2) | Out-Null
$d.Content.Find.Execute("87×14=1218", $true, $false, $false, $false, $false, $true, 1, $false, "26×53=1378", 2) | Out-Null
$d.Content.Find.Execute("72×62=4464", $true, $false, $false, $false, $false, $true, 1, $false, "70×65=4550", 2) | Out-Null
$d.Content.Find.Execute("15×88=1320", $true, $false, $false, $false, $false, $true, 1, $false, "92×99=9108", 2) | Out-Null
$d.Content.Find.Execute("22×55=1210", $true, $false, $false, $false, $false, $true, 1, $false, "36×73=2628", 2) | Out-Null
$d.Content.Find.Execute("23×85=1955", $true, $false, $false, $false, $false, $true, 1, $false, "61×55=3355", 2) | Out-Null
$d.Content.Find.Execute("13×60=780", $true, $false, $false, $false, $false, $true, 1, $false, "96×92=8832", 2) | Out-Null
$d.Content.Find.Execute("85×34=2890", $true, $false, $false, $false, $false, $true, 1, $false, "32×34=1088", 2) | Out-Null
$d.Content.Find.Execute("25×94=2350", $true, $false, $false, $false, $false, $true, 1, $false, "46×73=3358", 2) | Out-Null
$d.Content.Find.Execute("92×24=2208", $true, $false, $false, $false, $false, $true, 1, $false, "49×88=4312", 2) | Out-Null
$d.Content.Find.Execute("42×64=2688", $true, $false, $false, $false, $false, $true, 1, $false, "78×63=4914", 2) | Out-Null
$d.Content.Find.Execute("11×95=1045", $true, $false, $false, $false, $false, $true, 1, $false, "47×56=2632", 2) | Out-Null
$d.Content.Find.Execute("77×99=7623", $true, $false, $false, $false, $false, $true, 1, $false, "26×66=1716", 2) | Out-Null
$d.Content.Find.Execute("48×72=3456", $true, $false, $false, $false, $false, $true, 1, $false, "61×44=2684", 2) | Out-Null
$d.Content.Find.Execute("30×54=1620", $true, $false, $false, $false, $false, $true, 1, $false, "61×66=4026", 2) | Out-Null
$d.Content.Find.Execute("86×27=2322", $true, $false, $false, $false, $false, $true, 1, $false, "18×83=1494", 2) | Out-Null
$d.Content.Find.Execute("57×66=3762", $true, $false, $false, $false, $false, $true, 1, $false, "60×15=900", 2) | Out-Null
$d.Content.Find.Execute("83×35=2905", $true, $false, $false, $false, $false, $true, 1, $false, "15×90=1350", 2) | Out-Null
